$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 166. This shifts the existing
# rows 166-203 down to 167-204 (and the sheet's used range grows to R204).
$ws.Rows.Item(166).Insert()

# Populate the newly inserted row 166 with a new week of price data for
# "Betarraga" at "Macroferia Regional de Talca", matching the constant
# columns used by every other row in this sheet.
$ws.Cells.Item(166, 1).Value = 5
$ws.Cells.Item(166, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(166, 3).Value = "Maule"
$ws.Cells.Item(166, 4).Value = 44511
$ws.Cells.Item(166, 5).Value = 7
$ws.Cells.Item(166, 6).Value = 100114014
$ws.Cells.Item(166, 7).Value = "Betarraga"
$ws.Cells.Item(166, 8).Value = "Sin especificar"
$ws.Cells.Item(166, 9).Value = "Primera"
$ws.Cells.Item(166, 10).Value = 5000
$ws.Cells.Item(166, 11).Value = 500
$ws.Cells.Item(166, 12).Value = 500
$ws.Cells.Item(166, 13).Value = 500
$ws.Cells.Item(166, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(166, 15).Value = "Región del Maule"
$ws.Cells.Item(166, 16).Value = 100
$ws.Cells.Item(166, 17).Value = 5
$ws.Cells.Item(166, 18).Value = "Hortaliza"
